{"js": "// Replace each two-digit multiplication expression in the document body\n// with its new value, matching the exact (unique) original text for each.\nconst replacements = [\n  [\"34\u00d738=\", \"63\u00d756=\"],\n  [\"39\u00d749=\", \"36\u00d785=\"],\n  [\"74\u00d788=\", \"94\u00d749=\"],\n  [\"84\u00d742=\", \"47\u00d730=\"],\n  [\"47\u00d794=\", \"14\u00d752=\"],\n  [\"78\u00d753=\", \"26\u00d730=\"],\n  [\"60\u00d750=\", \"27\u00d773=\"],\n  [\"19\u00d739=\", \"33\u00d719=\"],\n  [\"38\u00d761=\", \"62\u00d769=\"],\n  [\"65\u00d772=\", \"97\u00d749=\"],\n  [\"82\u00d784=\", \"73\u00d777=\"],\n  [\"38\u00d783=\", \"79\u00d775=\"],\n  [\"29\u00d738=\", \"19\u00d722=\"],\n  [\"95\u00d791=\", \"67\u00d767=\"],\n  [\"74\u00d735=\", \"72\u00d712=\"],\n  [\"15\u00d796=\", \"93\u00d772=\"],\n  [\"96\u00d786=\", \"98\u00d775=\"],\n  [\"72\u00d728=\", \"50\u00d790=\"],\n  [\"48\u00d798=\", \"20\u00d747=\"],\n  [\"55\u00d778=\", \"95\u00d753=\"],\n  [\"31\u00d765=\", \"93\u00d750=\"],\n  [\"65\u00d760=\", \"75\u00d761=\"],\n  [\"80\u00d775=\", \"62\u00d751=\"],\n  [\"84\u00d773=\", \"55\u00d769=\"],\n  [\"22\u00d742=\", \"66\u00d757=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document\n# with its new value using Find/Replace (Execute) on the whole content.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"34\u00d738=\", \"63\u00d756=\"),\n  @(\"39\u00d749=\", \"36\u00d785=\"),\n  @(\"74\u00d788=\", \"94\u00d749=\"),\n  @(\"84\u00d742=\", \"47\u00d730=\"),\n  @(\"47\u00d794=\", \"14\u00d752=\"),\n  @(\"78\u00d753=\", \"26\u00d730=\"),\n  @(\"60\u00d750=\", \"27\u00d773=\"),\n  @(\"19\u00d739=\", \"33\u00d719=\"),\n  @(\"38\u00d761=\", \"62\u00d769=\"),\n  @(\"65\u00d772=\", \"97\u00d749=\"),\n  @(\"82\u00d784=\", \"73\u00d777=\"),\n  @(\"38\u00d783=\", \"79\u00d775=\"),\n  @(\"29\u00d738=\", \"19\u00d722=\"),\n  @(\"95\u00d791=\", \"67\u00d767=\"),\n  @(\"74\u00d735=\", \"72\u00d712=\"),\n  @(\"15\u00d796=\", \"93\u00d772=\"),\n  @(\"96\u00d786=\", \"98\u00d775=\"),\n  @(\"72\u00d728=\", \"50\u00d790=\"),\n  @(\"48\u00d798=\", \"20\u00d747=\"),\n  @(\"55\u00d778=\", \"95\u00d753=\"),\n  @(\"31\u00d765=\", \"93\u00d750=\"),\n  @(\"65\u00d760=\", \"75\u00d761=\"),\n  @(\"80\u00d775=\", \"62\u00d751=\"),\n  @(\"84\u00d773=\", \"55\u00d769=\"),\n  @(\"22\u00d742=\", \"66\u00d757=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Text = $old\n  $range.Find.Replacement.Text = $new\n  $range.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
